# Swap the identifying data between row 2 and row 4 (two observation records
# of "Knärot" at the same locality got their Id/Antal/Ålder-Stadium/
# Ost/Nord/Publik kommentar values exchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("A", "I", "K", "Q", "R", "AC")

foreach ($col in $columns) {
    $cellRow2 = $ws.Range($col + "2")
    $cellRow4 = $ws.Range($col + "4")

    $val2 = $cellRow2.Value2
    $val4 = $cellRow4.Value2

    $cellRow2.Value2 = $val4
    $cellRow4.Value2 = $val2
}
